$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.110.07'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.38%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.108.00'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.19%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '349.11'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.92%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.004'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.39%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5161'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.51%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4455'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.17%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '52.60'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -4.21%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08981'
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.176'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.43%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '25.80'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +5.09%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.111.59'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.30%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '8.301'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.27%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.753'
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '99.21'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.27%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001150'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.23%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.006'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.43%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '20.90'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +7.81%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.06687'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.06%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.004'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.35%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.246'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.36%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '30.221.57'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.24%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.87'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.79%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.347'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.28%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.359.63'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.21%  '
$ws.Range("E27").Value = '  -1.36%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.541'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.44%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '162.30'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.82%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '133.85'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.02%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.179'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.54%  '
$ws.Range("E32").Value = '  -0.21%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.635'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.45%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.261'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.61%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.967'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.36%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '10.41'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.19%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.942'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.66%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02579'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.47%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06836'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.37%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2310'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.58%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '12.62'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.41%  '
$ws.Range("E42").Value = '  -0.52%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.279'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.63%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.26'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.22%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.310'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.18%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6385'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.87%  '
$ws.Range("E47").Value = '  +2.55%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.658'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.77%  '
$ws.Range("E49").Value = '  -2.27%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '82.71'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.58%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07238'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.61%  '
